$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-11-12 Wednesday" "2025-11-13 Thursday"
Replace-Text "43-25=18" "47-4=43"
Replace-Text "61-37=24" "14+73=87"
Replace-Text "45-9=36" "99-89=10"
Replace-Text "46+26=72" "95-71=24"
Replace-Text "5+53=58" "4+37=41"
Replace-Text "90-87=3" "59+32=91"
Replace-Text "53-42=11" "8+38=46"
Replace-Text "90-11=79" "1+19=20"
Replace-Text "95-65=30" "47+13=60"
Replace-Text "90-66=24" "1+50=51"
Replace-Text "64+28=92" "63+7=70"
Replace-Text "12+42=54" "81-74=7"
Replace-Text "74-4=70" "78-50=28"
Replace-Text "56-40=16" "45+7=52"
Replace-Text "80-18=62" "56-41=15"
Replace-Text "52-6=46" "16+36=52"
Replace-Text "3+3=6" "65-38=27"
Replace-Text "14+68=82" "23-20=3"
Replace-Text "5+77=82" "88-56=32"
Replace-Text "57-18=39" "48-0=48"
Replace-Text "77-37=40" "87-11=76"
Replace-Text "94-7=87" "20+8=28"
Replace-Text "14-11=3" "81-43=38"
Replace-Text "8+23=31" "76-15=61"
Replace-Text "23+45=68" "43+24=67"
Replace-Text "48-39=9" "54+20=74"
Replace-Text "95-62=33" "27-10=17"
Replace-Text "9+35=44" "23+26=49"
Replace-Text "83-54=29" "36-18=18"
Replace-Text "71-56=15" "39-22=17"
Replace-Text "70+27=97" "75+16=91"
Replace-Text "32-16=16" "2+6=8"
Replace-Text "67+0=67" "19+23=42"
Replace-Text "3+72=75" "74-9=65"
Replace-Text "37-23=14" "1+47=48"
Replace-Text "66+2=68" "21+33=54"
Replace-Text "20+1=21" "70-4=66"
Replace-Text "54+6=60" "78-51=27"
Replace-Text "96-54=42" "42+12=54"
Replace-Text "0+78=78" "50-21=29"
Replace-Text "26+57=83" "72-51=21"
Replace-Text "51-9=42" "73+11=84"
Replace-Text "25+4=29" "10+15=25"
Replace-Text "31+6=37" "31+4=35"
Replace-Text "23+8=31" "36-31=5"
Replace-Text "44+7=51" "72-43=29"
Replace-Text "65-61=4" "55-39=16"
Replace-Text "33+57=90" "90-73=17"
Replace-Text "17+60=77" "66-36=30"
Replace-Text "30+43=73" "50+15=65"
Replace-Text "7+23=30" "92-18=74"
Replace-Text "52+20=72" "51-36=15"
Replace-Text "15+66=81" "5+31=36"
Replace-Text "78+15=93" "2+37=39"
Replace-Text "39-18=21" "50+13=63"
Replace-Text "3+33=36" "36+48=84"
Replace-Text "58-30=28" "72-37=35"
Replace-Text "55-40=15" "26-10=16"
Replace-Text "6+79=85" "77-38=39"
Replace-Text "73-31=42" "9-1=8"
Replace-Text "59-57=2" "23+75=98"
Replace-Text "9+66=75" "1+45=46"
Replace-Text "11+78=89" "95-82=13"
Replace-Text "21-12=9" "20+17=37"
Replace-Text "70-27=43" "31+24=55"
Replace-Text "59-50=9" "41+13=54"
Replace-Text "32+66=98" "46-12=34"
Replace-Text "73-21=52" "2+69=71"
Replace-Text "16+26=42" "99-77=22"
Replace-Text "42+24=66" "15+40=55"
Replace-Text "48+35=83" "92-27=65"
Replace-Text "4+31=35" "23+62=85"
Replace-Text "61+22=83" "76-19=57"
Replace-Text "77-54=23" "6+90=96"
Replace-Text "17+17=34" "69-36=33"
Replace-Text "29+53=82" "85-31=54"
Replace-Text "30+0=30" "1+71=72"
Replace-Text "78-2=76" "67+9=76"
Replace-Text "74-48=26" "73-51=22"
Replace-Text "28+62=90" "45+52=97"
Replace-Text "21-11=10" "31+14=45"
Replace-Text "95-79=16" "46-16=30"
Replace-Text "65-37=28" "90-4=86"
Replace-Text "34+49=83" "87+9=96"
Replace-Text "96-41=55" "61+31=92"
Replace-Text "5+41=46" "71-8=63"
Replace-Text "63-29=34" "43-24=19"
Replace-Text "0+50=50" "75+7=82"
Replace-Text "97-57=40" "70+6=76"
Replace-Text "99-5=94" "98-46=52"
Replace-Text "13-10=3" "50+48=98"
Replace-Text "78-46=32" "45+26=71"
Replace-Text "64-5=59" "20+22=42"
Replace-Text "44+45=89" "54+0=54"
Replace-Text "4+62=66" "89-48=41"
Replace-Text "73-30=43" "43+31=74"
Replace-Text "30-25=5" "64-7=57"
Replace-Text "8+35=43" "70-43=27"
Replace-Text "20+44=64" "45+13=58"
Replace-Text "20+74=94" "0+47=47"
